$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = $val
}

function Set-PercentTextValue($addr, $val, $donorAddr) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $ws.Range($donorAddr).Copy()
    $r.PasteSpecial(-4122)
}

Set-TextValue "E2" "2026-02-13 05:48:28"
Set-PercentTextValue "H2" "80%" "G2"
Set-TextValue "E3" "2026-02-13 05:48:30"
Set-TextValue "O3" "-2.6 °C"
Set-TextValue "E4" "2026-02-13 05:48:33"
Set-TextValue "J4" "1002.9 hPa"
Set-TextValue "O4" "9.8 °C"
Set-TextValue "E5" "2026-02-13 05:48:36"
Set-TextValue "L5" "36.7 km/h - 118º 5:22 TU"
Set-TextValue "E6" "2026-02-13 05:48:38"
Set-TextValue "J6" "1003.1 hPa"
Set-TextValue "E7" "2026-02-13 05:48:41"
Set-PercentTextValue "H7" "45%" "G7"
Set-TextValue "J7" "1003.6 hPa"
Set-TextValue "N7" "13.8 °C 5:09 TU"
Set-TextValue "O7" "14.5 °C"
Set-TextValue "E8" "2026-02-13 05:48:43"
Set-TextValue "J8" "1003.5 hPa"
Set-TextValue "N8" "9.8 °C 5:24 TU"
Set-TextValue "O8" "10.4 °C"
Set-TextValue "E9" "2026-02-13 05:48:45"
Set-PercentTextValue "H9" "65%" "G9"
Set-TextValue "O9" "8.5 °C"
Set-TextValue "E10" "2026-02-13 05:48:48"
Set-PercentTextValue "H10" "76%" "G10"
Set-TextValue "O10" "7.7 °C"
Set-TextValue "E11" "2026-02-13 05:48:50"
Set-TextValue "E12" "2026-02-13 05:48:53"
Set-TextValue "O12" "8.1 °C"
Set-TextValue "E13" "2026-02-13 05:48:55"
Set-PercentTextValue "H13" "84%" "G13"
Set-TextValue "J13" "1007.3 hPa"
Set-TextValue "O13" "-0.6 °C"
Set-TextValue "E14" "2026-02-13 05:48:57"
Set-PercentTextValue "H14" "61%" "G14"
Set-TextValue "N14" "9.0 °C 5:01 TU"
Set-TextValue "O14" "11.5 °C"
Set-TextValue "E15" "2026-02-13 05:49:00"
Set-PercentTextValue "H15" "66%" "G15"
Set-TextValue "O15" "8.8 °C"
Set-TextValue "E16" "2026-02-13 05:49:02"
Set-PercentTextValue "H16" "60%" "G16"
Set-TextValue "L16" "74.2 km/h - 285º 5:22 TU"
Set-TextValue "E17" "2026-02-13 05:49:05"
Set-TextValue "E18" "2026-02-13 05:49:07"
Set-TextValue "J18" "1003.2 hPa"
Set-TextValue "O18" "7.5 °C"
Set-TextValue "E19" "2026-02-13 05:49:10"
Set-PercentTextValue "H19" "70%" "G19"
Set-TextValue "E20" "2026-02-13 05:49:12"
Set-TextValue "E21" "2026-02-13 05:49:15"
Set-TextValue "J21" "1005.8 hPa"
Set-TextValue "E22" "2026-02-13 05:49:18"
Set-PercentTextValue "H22" "82%" "G22"
Set-TextValue "E23" "2026-02-13 05:49:20"
Set-TextValue "L23" "38.5 km/h - 98º 5:13 TU"
Set-TextValue "E24" "2026-02-13 05:49:23"
Set-PercentTextValue "H24" "85%" "G24"
Set-TextValue "I24" "0.1 mm"
Set-TextValue "J24" "1004.7 hPa"
Set-TextValue "E25" "2026-02-13 05:49:25"
Set-PercentTextValue "H25" "58%" "G25"
Set-TextValue "O25" "-2.7 °C"
Set-TextValue "E26" "2026-02-13 05:49:28"
Set-PercentTextValue "H26" "53%" "G26"
Set-TextValue "J26" "1003.6 hPa"
Set-TextValue "N26" "1.0 °C 5:26 TU"
Set-TextValue "O26" "2.7 °C"
Set-TextValue "E27" "2026-02-13 05:49:30"
Set-TextValue "O27" "-1.6 °C"
Set-TextValue "E28" "2026-02-13 05:49:33"
Set-TextValue "J28" "1003.7 hPa"
Set-TextValue "O28" "5.4 °C"
Set-TextValue "E29" "2026-02-13 05:49:35"
Set-TextValue "E30" "2026-02-13 05:49:38"
Set-TextValue "J30" "1003.4 hPa"
Set-TextValue "E31" "2026-02-13 05:49:40"
Set-PercentTextValue "H31" "53%" "G31"
Set-TextValue "J31" "1002.4 hPa"
Set-TextValue "N31" "10.7 °C 5:28 TU"
Set-TextValue "O31" "11.7 °C"
Set-TextValue "E32" "2026-02-13 05:49:43"
Set-TextValue "O32" "5.1 °C"
Set-TextValue "E33" "2026-02-13 05:49:45"
Set-PercentTextValue "H33" "75%" "G33"
Set-TextValue "J33" "1005.9 hPa"
Set-TextValue "O33" "0.7 °C"
Set-TextValue "E34" "2026-02-13 05:49:48"
Set-PercentTextValue "H34" "54%" "G34"
Set-TextValue "O34" "-0.3 °C"
Set-TextValue "E35" "2026-02-13 05:49:50"
Set-PercentTextValue "H35" "58%" "G35"
Set-TextValue "I35" "0.2 mm"
Set-TextValue "J35" "1005.4 hPa"
Set-TextValue "E36" "2026-02-13 05:49:53"
Set-PercentTextValue "H36" "59%" "G36"
Set-TextValue "J36" "1003.1 hPa"
Set-TextValue "N36" "8.1 °C 5:15 TU"
Set-TextValue "O36" "11.4 °C"
Set-TextValue "E37" "2026-02-13 05:49:55"
Set-PercentTextValue "H37" "63%" "G37"
Set-TextValue "J37" "1005.2 hPa"
Set-TextValue "O37" "3.9 °C"
Set-TextValue "E38" "2026-02-13 05:49:58"
Set-PercentTextValue "H38" "53%" "G38"
Set-TextValue "N38" "8.2 °C 5:03 TU"
Set-TextValue "O38" "10.5 °C"
Set-TextValue "E39" "2026-02-13 05:50:00"
Set-PercentTextValue "H39" "50%" "G39"
Set-TextValue "O39" "-2.5 °C"
Set-TextValue "E40" "2026-02-13 05:50:03"
Set-TextValue "J40" "1007.0 hPa"
Set-TextValue "E41" "2026-02-13 05:50:05"
Set-PercentTextValue "H41" "50%" "G41"
Set-TextValue "J41" "1004.1 hPa"
Set-TextValue "O41" "13.1 °C"
Set-TextValue "E42" "2026-02-13 05:50:08"
Set-TextValue "O42" "10.5 °C"
Set-TextValue "E43" "2026-02-13 05:50:10"
Set-PercentTextValue "H43" "65%" "G43"
Set-TextValue "N43" "4.9 °C 5:18 TU"
Set-TextValue "O43" "7.4 °C"
Set-TextValue "E44" "2026-02-13 05:50:12"
Set-PercentTextValue "H44" "80%" "G44"
Set-TextValue "L44" "59.0 km/h - 197º 5:10 TU"
Set-TextValue "E45" "2026-02-13 05:50:15"
Set-PercentTextValue "H45" "69%" "G45"
Set-TextValue "J45" "1004.2 hPa"
Set-TextValue "O45" "2.5 °C"
Set-TextValue "E46" "2026-02-13 05:50:17"
Set-PercentTextValue "H46" "81%" "G46"
Set-TextValue "J46" "1005.1 hPa"
